$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# Bump the date in A1 by one day (45308 -> 45309)
$ws.Range("A1").Value = Get-Date -Year 2024 -Month 1 -Day 18 -Hour 0 -Minute 0 -Second 0

# Step 1: VARILLA CHATA price update
$ws.Range("D22").Value = 288

# Step 2: ALAMBRE price update
$ws.Range("D38").Value = 367.127

# Refresh the merged ranges touched while editing steps 1 and 2
$ws.Range("B21:C21").UnMerge()
$ws.Range("B21:C21").Merge()
$ws.Range("B22:C22").UnMerge()
$ws.Range("B22:C22").Merge()
